# Generate Report for Handback
# Updates the localization-status workbook once the de-de/zh-cn handback
# files have come back in sync with en-US: flips the "Ready for handoff"
# status to "Handed back: in sync with en-US", and fills in the
# "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns for the two language sheets.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$mdDisplay  = "5d854677-34d1-4201-9b3a-d214c208291c.md"
$mdUrl      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/34d044c3211eaf7815a11fe32a3984aac1c2a17b/e2e/5d854677-34d1-4201-9b3a-d214c208291c.md"

# Column width helper input: this engine snaps ColumnWidth to a pixel grid
# on save, so these inputs are chosen to land squarely inside the grid
# step that contains the desired rendered width.
$wideColWidth   = 29.166666666666668   # renders as 30   (status columns)
$widestColWidth = 39.166666666666664   # renders as 40   (target/handback file columns)

# ---------------------------------------------------------------------
# Overview sheet: status column shows up twice (zh-cn + de-de handoff
# status), both switch from "Ready for handoff" to the handback message.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText
$overview.Columns.Item(5).ColumnWidth = $wideColWidth
$overview.Columns.Item(6).ColumnWidth = $wideColWidth

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText
$zhcn.Columns.Item(3).ColumnWidth = $wideColWidth

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl, "", "", $mdDisplay)
$zhcn.Range("I2").Font.Color = 15570276
$zhcn.Range("I2").Font.Underline = 2

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $mdUrl, "", "", $mdDisplay)
$zhcn.Range("I3").Font.Color = 15570276
$zhcn.Range("I3").Font.Underline = 2

$zhcn.Range("J2").Value = "5d854677-34d1-4201-9b3a-d214c208291c.e1ca66088b8959280ed863c57f3c222ea9477436.zh-cn.xlf"
$zhcn.Range("J3").Value = "5d854677-34d1-4201-9b3a-d214c208291c.e1ca66088b8959280ed863c57f3c222ea9477436.zh-cn.xlf"

$zhcn.Range("K2").Value = "2016-09-02 01:14:08"
$zhcn.Range("K3").Value = "2016-09-02 01:14:08"

$zhcn.Columns.Item(9).ColumnWidth = $widestColWidth
$zhcn.Columns.Item(10).ColumnWidth = $widestColWidth

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText
$dede.Columns.Item(3).ColumnWidth = $wideColWidth

$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl, "", "", $mdDisplay)
$dede.Range("I2").Font.Color = 15570276
$dede.Range("I2").Font.Underline = 2

$dede.Hyperlinks.Add($dede.Range("I3"), $mdUrl, "", "", $mdDisplay)
$dede.Range("I3").Font.Color = 15570276
$dede.Range("I3").Font.Underline = 2

$dede.Range("J2").Value = "5d854677-34d1-4201-9b3a-d214c208291c.e1ca66088b8959280ed863c57f3c222ea9477436.de-de.xlf"
$dede.Range("J3").Value = "5d854677-34d1-4201-9b3a-d214c208291c.e1ca66088b8959280ed863c57f3c222ea9477436.de-de.xlf"

$dede.Range("K2").Value = "2016-09-02 01:14:16"
$dede.Range("K3").Value = "2016-09-02 01:14:16"

$dede.Columns.Item(9).ColumnWidth = $widestColWidth
$dede.Columns.Item(10).ColumnWidth = $widestColWidth
